$d = $word.ActiveDocument

# --- 1) "Explorer, Google Chrome, Safari" -- merge the split runs (remove the
#        proofErr spell-check markers around "Chrome") into a single run.
$d.Content.Find.Execute(
    "Explorer, Google Chrome, Safari", $false, $false, $false, $false, $false,
    $true, 1, $false, "Explorer, Google Chrome, Safari", 2) | Out-Null

# --- 2) "El sistema mobile debe contar..." -- merge the split runs (remove the
#        proofErr spell-check markers around "mobile") into a single run.
$d.Content.Find.Execute(
    "El sistema mobile debe contar con textos que tengan un lenguaje amigable al usuario evitando el uso de términos técnicos",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "El sistema mobile debe contar con textos que tengan un lenguaje amigable al usuario evitando el uso de términos técnicos",
    2) | Out-Null

# --- 3) "...de forma encriptada" -- merge the split runs (remove the proofErr
#        spell-check markers around "encriptada") into a single run.
$d.Content.Find.Execute(
    "Las contraseñas se almacenaran en base de datos de forma encriptada",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Las contraseñas se almacenaran en base de datos de forma encriptada",
    2) | Out-Null

# --- 4) "Únicamente ... de administrados podrán ..." -- fix the typo
#        ("administrados" -> "administrador"), then split the run so that
#        "de administrador" is its own run, and move the "_GoBack" bookmark
#        to sit right after it (this is what happened in real Word: the
#        _GoBack bookmark tracks the last edited location).
$d.Content.Find.Execute(
    "de administrados", $false, $false, $false, $false, $false,
    $true, 1, $false, "de administrador", 2) | Out-Null

$adminPara = $d.Paragraphs.Item(14)
$adminRange = $adminPara.Range.Duplicate
$adminFind = $adminRange.Find
$adminFind.ClearFormatting()
$adminFind.Text = "de administrador"
$adminFind.Execute() | Out-Null
# Touching a character-formatting property on the sub-range forces Word to
# split it into its own run without altering the final formatting.
$adminRange.Bold = 1
$adminRange.Bold = 0

$goBackPoint = $d.Range($adminRange.End, $adminRange.End)
$d.Bookmarks.Add("_GoBack", $goBackPoint) | Out-Null

# --- 5) "El sistema  web debe contar con un sistema sencillo de deploy en el
#        servidor" -- merge the split runs (remove the proofErr spell-check
#        markers around "deploy") into a single run. (Note: double space
#        between "sistema" and "web" is preserved.)
$d.Content.Find.Execute(
    "El sistema  web debe contar con un sistema sencillo de deploy en el servidor",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "El sistema  web debe contar con un sistema sencillo de deploy en el servidor",
    2) | Out-Null

# --- 6) "El proceso de desarrollo ... deployarlos en producción ..." -- merge
#        the split runs (remove the proofErr spell-check markers around
#        "deployarlos") into a single run.
$d.Content.Find.Execute(
    "El proceso de desarrollo debe permitir subir cambios al servidor web y deployarlos en producción de forma sencilla para minimizar tiempo de espera de resolución de fallos y de indisponibilidad del servicio",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "El proceso de desarrollo debe permitir subir cambios al servidor web y deployarlos en producción de forma sencilla para minimizar tiempo de espera de resolución de fallos y de indisponibilidad del servicio",
    2) | Out-Null

# --- 7) "Los sistemas tanto mobile como web ..." -- merge the split runs
#        (remove the proofErr spell-check markers around "mobile") into a
#        single run; "homogénea" remains its own trailing run.
$d.Content.Find.Execute(
    "Los sistemas tanto mobile como web deben contar con interfaces de usuarios definidas de forma ",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Los sistemas tanto mobile como web deben contar con interfaces de usuarios definidas de forma ",
    2) | Out-Null
